$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the existing AutoFilter criterion on column W ("Management") so
# --- every previously-filtered-out row becomes visible again.
$ws.ShowAllData()

# --- New "Change in food prices" mini-table (AB1:AE5) ---
$ws.Range("AB1").Value = 'Energy sector changes'
$ws.Range("AC1").Value = 2050
$ws.Range("AD1").Value = 'SSP'
$ws.Range("AE1").Value = 'Units'
$ws.Range("AB2").Value = 'Non-Staples'
$ws.Range("AC2").Value = 0.034622
$ws.Range("AD2").Value = 'SSP1'
$ws.Range("AE2").Value = 'USD$/Mcal/day'
$ws.Range("AB3").Value = 'Staples'
$ws.Range("AC3").Value = -0.002542
$ws.Range("AD3").Value = 'SSP1'
$ws.Range("AE3").Value = 'USD$/Mcal/day'
$ws.Range("AB4").Value = 'Non-Staples'
$ws.Range("AC4").Value = 13.283339
$ws.Range("AD4").Value = 'SSP1'
$ws.Range("AE4").Value = '%'
$ws.Range("AB5").Value = 'Staples'
$ws.Range("AC5").Value = -2.760822
$ws.Range("AD5").Value = 'SSP1'
$ws.Range("AE5").Value = '%'

# --- "Food Sector Carbon Emissions in 2050" block (A21:C30) ---
$ws.Range("A21").Value = 'Food Sector Carbon Emissions in 2050'
$ws.Range("A22").Value = 'GWP of methane: '
$ws.Range("B22").Value = 27.2
$ws.Range("A23").Value = 'CO2-Ceq'
$ws.Range("B23").Formula = "=12/44"
$ws.Range("A24").Value = 'Methane C avoidance'
$ws.Range("B24").Formula = "=H18*B22*B23"
$ws.Range("A25").Value = 'Carbon Sequestration'
$ws.Range("B25").Formula = "=B18"
$ws.Range("A26").Value = 'Total'
$ws.Range("B26").Formula = "=B24+B25"
$ws.Range("A28").Value = 'Food system emissions in 2050 (Gt CO2-eq)'
$ws.Range("B28").Value = 18.4
$ws.Range("C28").Value = '[Almaraz et al., 2023]'
$ws.Range("A29").Value = 'Food system emissions in 2050 (Mt C)'
$ws.Range("B29").Formula = "=B28*1000*B23"
$ws.Range("A30").Value = 'Percentage reduction'
$ws.Range("B30").Formula = "=100*B26/B29"

# --- "Energy sector changes" table #1: absolute EJ changes (G21:J31) ---
$ws.Range("G21").Value = 'Energy sector changes'
$ws.Range("H21").Value = 2050
$ws.Range("I21").Value = 'fuel'
$ws.Range("J21").Value = 'Units'

$ws.Range("G22").Value = 0
$ws.Range("H22").Value = -5.449188
$ws.Range("I22").Value = 'oil'
$ws.Range("J22").Value = 'EJ'

$ws.Range("G23").Value = 1
$ws.Range("H23").Value = -13.864243
$ws.Range("I23").Value = 'natural gas'
$ws.Range("J23").Value = 'EJ'

$ws.Range("G24").Value = 2
$ws.Range("H24").Value = -15.926671
$ws.Range("I24").Value = 'coal'
$ws.Range("J24").Value = 'EJ'

$ws.Range("G25").Value = 3
$ws.Range("H25").Value = -5.265949
$ws.Range("I25").Value = 'biomass'
$ws.Range("J25").Value = 'EJ'

$ws.Range("G26").Value = 4
$ws.Range("H26").Value = 5.105476
$ws.Range("I26").Value = 'nuclear'
$ws.Range("J26").Value = 'EJ'

$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 'hydro'
$ws.Range("J27").Value = 'EJ'

$ws.Range("G28").Value = 6
$ws.Range("H28").Value = 5.704001
$ws.Range("I28").Value = 'wind'
$ws.Range("J28").Value = 'EJ'

$ws.Range("G29").Value = 7
$ws.Range("H29").Value = 5.833608
$ws.Range("I29").Value = 'solar'
$ws.Range("J29").Value = 'EJ'

$ws.Range("G30").Value = 8
$ws.Range("H30").Value = 0.596756
$ws.Range("I30").Value = 'geothermal'
$ws.Range("J30").Value = 'EJ'

$ws.Range("G31").Value = 9
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 'traditional biomass'
$ws.Range("J31").Value = 'EJ'

# --- "Energy sector changes" table #2: percentage changes (H32:J42) ---
$ws.Range("H32").Value = 2050
$ws.Range("I32").Value = 'fuel'
$ws.Range("J32").Value = 'Units'

$ws.Range("G33").Value = 0
$ws.Range("H33").Value = -3.349711
$ws.Range("I33").Value = 'oil'
$ws.Range("J33").Value = '%'

$ws.Range("G34").Value = 1
$ws.Range("H34").Value = -8.343863
$ws.Range("I34").Value = 'natural gas'
$ws.Range("J34").Value = '%'

$ws.Range("G35").Value = 2
$ws.Range("H35").Value = -25.94169
$ws.Range("I35").Value = 'coal'
$ws.Range("J35").Value = '%'

$ws.Range("G36").Value = 3
$ws.Range("H36").Value = -2.905689
$ws.Range("I36").Value = 'biomass'
$ws.Range("J36").Value = '%'

$ws.Range("G37").Value = 4
$ws.Range("H37").Value = 15.637447
$ws.Range("I37").Value = 'nuclear'
$ws.Range("J37").Value = '%'

$ws.Range("G38").Value = 5
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 'hydro'
$ws.Range("J38").Value = '%'

$ws.Range("G39").Value = 6
$ws.Range("H39").Value = 7.50808
$ws.Range("I39").Value = 'wind'
$ws.Range("J39").Value = '%'

$ws.Range("G40").Value = 7
$ws.Range("H40").Value = 6.371349
$ws.Range("I40").Value = 'solar'
$ws.Range("J40").Value = '%'

$ws.Range("G41").Value = 8
$ws.Range("H41").Value = 7.759209
$ws.Range("I41").Value = 'geothermal'
$ws.Range("J41").Value = '%'

$ws.Range("G42").Value = 9
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 'traditional biomass'
$ws.Range("J42").Value = '%'

# --- Restore the selection / scroll state Excel saved after the edit ---
$ws.Range("G22").Select()
